# Feria Lagunitas de Puerto Montt - Piña: insert a new weekly price record
# as row 65, pushing the existing rows 65..180 down to 66..181.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 65 (shifts rows 65:180 down to 66:181)
$ws.Rows("65:65").Insert()

# Populate the newly inserted row 65 with the new record's data
$ws.Range("A65").Value = 4
$ws.Range("B65").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C65").Value = "Los Lagos"
$ws.Range("D65").Value = 44557
$ws.Range("E65").Value = 10
$ws.Range("F65").Value = "Fruta"
$ws.Range("G65").Value = 100108
$ws.Range("H65").Value = "Tropicales y subtropicales"
$ws.Range("I65").Value = 100108005
$ws.Range("J65").Value = "Piña"
$ws.Range("K65").Value = "Caramelo"
$ws.Range("L65").Value = "Tercera"
$ws.Range("M65").Value = 120
$ws.Range("N65").Value = 19000
$ws.Range("O65").Value = 20000
$ws.Range("P65").Value = 19500
$ws.Range("Q65").Value = "$/caja 16 unidades"
$ws.Range("R65").Value = "Ecuador"
$ws.Range("S65").Value = 1219
$ws.Range("T65").Value = 16
